# nestjs-03 -> session 04: add exercise bullets to the "Excercise" slide
# (the slide titled "Excercise", containing the "Required:"/"Optional:" list).
#
# Target structure for the Content Placeholder text body:
#   Required:
#       Exercise all the stuff we discuss in this session yourself   (lvl 1)
#   Optional:
#       How to avoid a use from liking a post twice.                 (lvl 1)
#       Complete other relations and entities                        (lvl 1)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# Locate the body placeholder. Prefer the well-known shape name, but fall
# back to matching on its current text in case the name ever changes.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $candidate = $s.Shapes.Item($i)
        if ($candidate.TextFrame.TextRange.Text -eq "Required:`rOptional:") {
            $shape = $candidate
        }
    }
}

$tr = $shape.TextFrame.TextRange

# New paragraph texts, in the final order.
$line1 = "Required:"
$line2 = "Exercise all the stuff we discuss in this session yourself"
$line3 = "Optional:"
$line4 = "How to avoid a use from liking a post twice."
$line5 = "Complete other relations and entities"

# Paragraph marks are carriage returns, just like real PowerPoint text ranges.
$tr.Text = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4 + "`r" + $line5

# Compute 1-based character offsets of the new sub-bullets so we can demote
# just those three paragraphs to outline level 2 (OOXML lvl="1") without
# touching the "Required:"/"Optional:" header paragraphs.
$pos2 = $line1.Length + 2
$pos4 = $pos2 + $line2.Length + 1 + $line3.Length + 1
$pos5 = $pos4 + $line4.Length + 1

$tr.Characters($pos2, $line2.Length).IndentLevel = 2
$tr.Characters($pos4, $line4.Length).IndentLevel = 2
$tr.Characters($pos5, $line5.Length).IndentLevel = 2
